# Update "Pais" (countries) worksheet with the newer snapshot of data.
# - Refresh the "last updated" timestamp label in A1.
# - Update case counters for several countries (rows shifted because some
#   countries' updated totals changed their relative ranking, causing the
#   country name in that row to change along with the new figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp label
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 07:56"

# Catar / Israel (rows 29-30) - Israel's updated numbers overtake Catar
$ws.Range("A29").Value = "Israel"
$ws.Range("B29").Value = 119265
$ws.Range("C29").Value = 727
$ws.Range("D29").Value = 96654
$ws.Range("E29").Value = 21654
$ws.Range("H29").Value = 957

$ws.Range("A30").Value = "Catar"
$ws.Range("B30").Value = 118994
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 115895
$ws.Range("E30").Value = 2901
$ws.Range("H30").Value = 198

# Kirguistan (row 59) - numbers refreshed, no row reordering
$ws.Range("B59").Value = 44036
$ws.Range("C59").Value = 78
$ws.Range("D59").Value = 38895
$ws.Range("E59").Value = 4082

# Costa Rica / Uzbekistan (rows 62-63) - Uzbekistan's updated numbers overtake Costa Rica
$ws.Range("A62").Value = "Uzbekistan"
$ws.Range("B62").Value = 42370
$ws.Range("C62").Value = 243
$ws.Range("D62").Value = 39664
$ws.Range("E62").Value = 2380
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 326

$ws.Range("A63").Value = "Costa Rica"
$ws.Range("B63").Value = 42184
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 16270
$ws.Range("E63").Value = 25471
$ws.Range("H63").Value = 443

# Afganistan (row 65) - numbers refreshed, no row reordering
$ws.Range("B65").Value = 38205
$ws.Range("C65").Value = 9
$ws.Range("D65").Value = 29254
$ws.Range("E65").Value = 7545

# Australia (row 72) - numbers refreshed, no row reordering
$ws.Range("B72").Value = 25923
$ws.Range("C72").Value = 104
$ws.Range("D72").Value = 21690
$ws.Range("E72").Value = 3570

# Sudan del Sur / Jamaica (rows 135-136) - Jamaica's updated numbers overtake Sudan del Sur
$ws.Range("A135").Value = "Jamaica"
$ws.Range("B135").Value = 2683
$ws.Range("C135").Value = 224
$ws.Range("D135").Value = 898
$ws.Range("E135").Value = 1761
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 24

$ws.Range("A136").Value = "Sudan del Sur"
$ws.Range("B136").Value = 2527
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 1290
$ws.Range("E136").Value = 1190
$ws.Range("H136").Value = 47
